# This script emulates a Power Query "Refresh" of the 쿼리2 query table:
# the monthly accumulated donation totals (column C) were re-pulled from the
# live API and the "새로고침시간" (refresh time) column (D) was stamped with
# the new refresh timestamp. The active selection also moved (a side effect
# of the user clicking cell F4 before saving), and column D's displayed
# width grew slightly to fit the new timestamp text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "월별 누적별풍선" (monthly cumulative total) values from the refresh
$ws.Range("C2").Value = 731973
$ws.Range("C4").Value = 533126
$ws.Range("C5").Value = 519289
$ws.Range("C6").Value = 428439
$ws.Range("C7").Value = 383363
$ws.Range("C9").Value = 278998
$ws.Range("C10").Value = 200078
$ws.Range("C11").Value = 152599

# New refresh timestamp written into "새로고침시간" for every row (the query
# only computes RefreshTime once, then stamps it onto every row it returns).
# Written as the raw OLE Automation date serial so sub-second precision
# survives (assigning a .NET DateTime truncates to whole seconds).
$ws.Range("D2:D12").Value = 46015.9792794213

# Column D got a hair wider to fit the refreshed text
$ws.Columns("D").ColumnWidth = 18.85

# Selection moved to F4 before the file was saved
$ws.Range("F4").Select()
